# Add LaborProductivity2015 column (S) to sheet1, mirroring the existing
# LaborProductivity2012 column (R).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header
$ws.Range("S1").Value = "LaborProductivity2015"
$ws.Range("S1").Style = $ws.Range("R1").Style

# Set column width to match the diff (target OOXML width 29.5703125;
# ColumnWidth is quantized to whole pixels at save time, so 28.6666... is
# the closest setting that round-trips to 29.5).
$ws.Columns.Item(19).ColumnWidth = 28.6666666666667

# Data values for S2:S18
$values = @{
    2  = 583368
    3  = 117327
    4  = 175870
    5  = 88532
    6  = 170684
    7  = 256386
    8  = 94144
    9  = 65196
    10 = 95211
    11 = 150228
    12 = 150555
    13 = 113384
    14 = 142402
    15 = 157780
    16 = 115772
    17 = 88178
    18 = 40008
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 19)
    $cell.Value = $values[$row]
    if ($row -eq 2) {
        $cell.Style = $ws.Range("R2").Style
    } else {
        $cell.NumberFormat = "#,##0"
        $cell.HorizontalAlignment = -4108  # xlCenter
    }
}

# Update the selection to reflect the new last column (S), keeping the
# existing freeze (column A) in place.
$ws.Range("S1").Select() | Out-Null
